$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $Address, $Text)
    $cell = $Worksheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "67.820.13"
Set-TextValue $ws "E2" "  +0.69%  "
Set-TextValue $ws "D3" "2.619.62"
Set-TextValue $ws "E3" "  -0.34%  "
Set-TextValue $ws "E4" "  -0.08%  "
Set-TextValue $ws "D5" "595.99"
Set-TextValue $ws "E5" "  -1.02%  "
Set-TextValue $ws "D6" "152.66"
Set-TextValue $ws "E6" "  -0.51%  "
Set-TextValue $ws "E7" "  -0.04%  "
Set-TextValue $ws "D8" "0.544"
Set-TextValue $ws "E8" "  -2.30%  "
Set-TextValue $ws "D9" "2.618.32"
Set-TextValue $ws "E9" "  -0.29%  "
Set-TextValue $ws "D10" "0.132"
Set-TextValue $ws "E10" "  +6.90%  "
Set-TextValue $ws "E11" "  -0.71%  "
Set-TextValue $ws "D12" "5.20"
Set-TextValue $ws "E12" "  -0.08%  "
Set-TextValue $ws "E13" "  -1.46%  "
Set-TextValue $ws "D14" "27.55"
Set-TextValue $ws "E14" "  -1.29%  "
Set-TextValue $ws "E15" "  +2.85%  "
Set-TextValue $ws "D16" "3.088.07"
Set-TextValue $ws "E16" "  -0.82%  "
Set-TextValue $ws "D17" "67.704.90"
Set-TextValue $ws "E17" "  +0.58%  "
Set-TextValue $ws "D18" "2.615.68"
Set-TextValue $ws "E18" "  -0.52%  "
Set-TextValue $ws "D19" "372.47"
Set-TextValue $ws "E19" "  +2.36%  "
Set-TextValue $ws "D20" "11.21"
Set-TextValue $ws "E20" "  -0.19%  "
Set-TextValue $ws "E21" "  -1.74%  "
Set-TextValue $ws "E22" "  -13.07%  "
Set-TextValue $ws "E23" "  -2.97%  "
Set-TextValue $ws "E24" "  -4.09%  "
Set-TextValue $ws "D25" "72.96"
Set-TextValue $ws "E25" "  +8.16%  "
Set-TextValue $ws "D27" "9.86"
Set-TextValue $ws "E27" "  -2.50%  "
Set-TextValue $ws "D28" "592.42"
Set-TextValue $ws "E28" "  +2.28%  "
Set-TextValue $ws "E30" "  +0.09%  "
Set-TextValue $ws "D31" "0.999"
Set-TextValue $ws "E31" "  -0.23%  "
Set-TextValue $ws "E32" "  -1.88%  "
Set-TextValue $ws "E33" "  -1.13%  "
Set-TextValue $ws "E34" "  -0.49%  "
Set-TextValue $ws "D35" "0.999"
Set-TextValue $ws "E35" "  -0.02%  "
Set-TextValue $ws "E36" "  -1.73%  "
Set-TextValue $ws "E37" "  -1.37%  "
Set-TextValue $ws "D38" "158.35"
Set-TextValue $ws "E38" "  +0.03%  "
Set-TextValue $ws "E39" "  -1.33%  "
Set-TextValue $ws "E40" "  +3.41%  "
Set-TextValue $ws "E41" "  -0.79%  "
Set-TextValue $ws "D42" "5.28"
Set-TextValue $ws "E42" "  -0.13%  "
Set-TextValue $ws "D43" "2.70"
Set-TextValue $ws "E43" "  +2.94%  "
Set-TextValue $ws "D44" "17.12"
Set-TextValue $ws "E44" "  +4.65%  "
Set-TextValue $ws "D45" "1.00"
Set-TextValue $ws "E45" "  +0.05%  "
Set-TextValue $ws "E46" "  -2.02%  "
Set-TextValue $ws "D47" "0.0₆0302"
Set-TextValue $ws "E47" "  +5.03%  "
Set-TextValue $ws "D48" "156.53"
Set-TextValue $ws "E48" "  +0.28%  "
Set-TextValue $ws "E49" "  -1.39%  "
Set-TextValue $ws "E50" "  -1.97%  "
Set-TextValue $ws "E51" "  -1.66%  "
